$d = $word.ActiveDocument

# --- Step 1: drop the two orphan bookmark markers that wrap paragraph 1 ---
# Deleting paragraph 1's whole range (including its paragraph mark) collapses
# both markers to the very start of the document, adjacent to one another.
$p1 = $d.Paragraphs.Item(1)
$d.Range($p1.Range.Start, $p1.Range.End).Delete()
# Consume the two now-adjacent zero-width markers one at a time.
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# --- Step 2: insert the new "Title" style heading paragraph, split word-by-word ---
$titleXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Articles</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">on</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Distributism</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Range(0, 0).InsertXML($titleXml)

# --- Step 3: turn the old "By Dorothy Day" paragraph into the "Authors" paragraph ---
$p2 = $d.Paragraphs.Item(2)
$authorsXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Range($p2.Range.Start, $p2.Range.End).InsertXML($authorsXml)
$d.Paragraphs.Item(2).Style = "Authors"
